$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 09:42"

# Row 4: Estados Unidos - updated case counts
$ws.Range("B4").Value = 2297338
$ws.Range("C4").Value = 148
$ws.Range("D4").Value = 956077
$ws.Range("E4").Value = 1219854

# Row 34: Singapur - updated case counts
$ws.Range("B34").Value = 41833
$ws.Range("C34").Value = 218
$ws.Range("E34").Value = 8348

# Row 38: Ucrania - updated case counts
$ws.Range("B38").Value = 35825
$ws.Range("C38").Value = 841
$ws.Range("D38").Value = 16406
$ws.Range("E38").Value = 18425
$ws.Range("G38").Value = 9
$ws.Range("H38").Value = 994

# Row 52: Armenia - updated case counts
$ws.Range("B52").Value = 19708
$ws.Range("C52").Value = 551
$ws.Range("D52").Value = 8854
$ws.Range("E52").Value = 10522
$ws.Range("G52").Value = 13
$ws.Range("H52").Value = 332

# Row 84: Gabon -> El Salvador
$ws.Range("A84").Value = "El Salvador"
$ws.Range("B84").Value = 4475
$ws.Range("C84").Value = 146
$ws.Range("D84").Value = 2449
$ws.Range("E84").Value = 1933
$ws.Range("G84").Value = 7
$ws.Range("H84").Value = 93

# Row 85: Kenia -> Gabon
$ws.Range("A85").Value = "Gabon"
$ws.Range("B85").Value = 4428
$ws.Range("D85").Value = 1750
$ws.Range("E85").Value = 2644
$ws.Range("H85").Value = 34

# Row 86: El Salvador -> Kenia
$ws.Range("A86").Value = "Kenia"
$ws.Range("B86").Value = 4374
$ws.Range("D86").Value = 1550
$ws.Range("E86").Value = 2705
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 119

# Row 88: Hungria - updated case counts
$ws.Range("B88").Value = 4086
$ws.Range("C88").Value = 5
$ws.Range("D88").Value = 2585
$ws.Range("E88").Value = 931
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 570

# Row 104: Estonia - updated case counts
$ws.Range("B104").Value = 1981
$ws.Range("C104").Value = 2
$ws.Range("D104").Value = 1758
$ws.Range("E104").Value = 154

# Row 114: Eslovaquia - updated case counts
$ws.Range("B114").Value = 1586
$ws.Range("C114").Value = 10
$ws.Range("E114").Value = 111

# Row 125: Letonia - updated case counts
$ws.Range("B125").Value = 1111
$ws.Range("C125").Value = 1
$ws.Range("E125").Value = 178

# Row 139: San Marino -> Estado de Palestina
$ws.Range("A139").Value = "Estado de Palestina"
$ws.Range("B139").Value = 707
$ws.Range("C139").Value = 32
$ws.Range("D139").Value = 437
$ws.Range("E139").Value = 267
$ws.Range("H139").Value = 3

# Row 140: Santo Tome y Principe -> San Marino
$ws.Range("A140").Value = "San Marino"
$ws.Range("B140").Value = 696
$ws.Range("D140").Value = 610
$ws.Range("E140").Value = 44
$ws.Range("H140").Value = 42

# Row 141: Estado de Palestina -> Santo Tome y Principe
$ws.Range("A141").Value = "Santo Tome y Principe"
$ws.Range("B141").Value = 693
$ws.Range("D141").Value = 199
$ws.Range("E141").Value = 482
$ws.Range("H141").Value = 12

# Row 202: Dominica -> Fiyi
$ws.Range("A202").Value = "Fiyi"

# Row 203: Fiyi -> Dominica
$ws.Range("A203").Value = "Dominica"

# Row 208: Islas Turcas y Caicos -> Santa Sede
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209: Santa Sede -> Islas Turcas y Caicos
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

